$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04011141694006874
$ws.Range("H2").Value = -16.83412728487585
$ws.Range("I2").Value = 50.47713347931106
$ws.Range("G3").Value = 0.0627292308392918
$ws.Range("H3").Value = 63.54219394819556
$ws.Range("G4").Value = -0.4707463410059471
$ws.Range("H4").Value = -3.970764937719685
$ws.Range("G5").Value = -0.4755622788970378
$ws.Range("H5").Value = 0.6834713178894997
$ws.Range("G6").Value = 0.2472773279702291
$ws.Range("H6").Value = 5.840680022377611
$ws.Range("G7").Value = 0.2626512616359919
$ws.Range("H7").Value = 19.07592614795621
$ws.Range("G8").Value = 0.1652813337650478
$ws.Range("H8").Value = -0.9142958902711003
$ws.Range("G9").Value = 0.1755788339174719
$ws.Range("H9").Value = 2.073993711795682
$ws.Range("G10").Value = -0.0006829463226870007
$ws.Range("H10").Value = 85.59027452859851
$ws.Range("G11").Value = 0.002772937684923746
$ws.Range("H11").Value = 118.9210627126073
$ws.Range("G12").Value = 0.1299940665913374
$ws.Range("H12").Value = -4.921076970905529
$ws.Range("G13").Value = 0.1444417545757515
$ws.Range("H13").Value = 15.8894067797217
$ws.Range("G14").Value = 0.2543523514079327
$ws.Range("H14").Value = 2.837818772357681
$ws.Range("G15").Value = 0.2634471138946617
$ws.Range("H15").Value = 4.262399632287157
$ws.Range("G16").Value = 0.1437213248549029
$ws.Range("H16").Value = -6.351725028143282
$ws.Range("G17").Value = 0.1492918236684211
$ws.Range("H17").Value = -1.141192908722041
$ws.Range("G18").Value = -0.0004048116894381696
$ws.Range("H18").Value = 97.52785339265166
$ws.Range("G19").Value = 0.007545671958485558
$ws.Range("H19").Value = 996.0925460823316
$ws.Range("G20").Value = 0.1348772251343799
$ws.Range("H20").Value = -2.729618090906378
$ws.Range("G21").Value = 0.1511009595328182
$ws.Range("H21").Value = 5.595243513432582
$ws.Range("G22").Value = 0.1592772411451793
$ws.Range("H22").Value = -14.46704060047433
$ws.Range("G23").Value = 0.1824330351587339
$ws.Range("H23").Value = 1.656370544943635
$ws.Range("G24").Value = -0.09208699941680973
$ws.Range("H24").Value = 2.447383751544725
$ws.Range("G25").Value = -0.07959933242092411
$ws.Range("H25").Value = 20.08777491271972
$ws.Range("G26").Value = 0.2246061005965345
$ws.Range("H26").Value = -2.404767263703069
$ws.Range("G27").Value = 0.2438368020071973
$ws.Range("H27").Value = 4.843468788889502
$ws.Range("G28").Value = 0.06245067125083205
$ws.Range("H28").Value = 6.202731632429629
$ws.Range("G29").Value = 0.08041253103060982
$ws.Range("H29").Value = 13.92374620370028
